$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes & Metrics")

# ---- Column width update (Sample Values column E widened) ----
$ws.Columns.Item(5).ColumnWidth = 39.85546875

# ---- Row 2 header restyle (fill4 stays, now center/center/wrap for most, center/center for C2) ----
$hdr = $ws.Range("A2:E2")
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4108    # xlCenter
$hdr.WrapText = $true
$ws.Range("C2").WrapText = $false
$ws.Rows.Item(2).RowHeight = 30

# ---- Rows 3-7: Album dimension block ----
# Merge A3:A7 for "Album"
$albumBlock = $ws.Range("A3:A7")
$albumBlock.Merge()
$ws.Range("A3").Value = "Album"

# Whole block A3:E7 gets fill D6DCE4, thin theme1 border, left/top/wrap
$block = $ws.Range("A3:E7")
$block.Interior.Pattern = 1
$block.Interior.Color = 14998742   # BGR encoding of RGB D6DCE4
$block.HorizontalAlignment = -4131  # xlLeft
$block.VerticalAlignment = -4160    # xlTop
$block.WrapText = $true

# Borders: thin, theme1(black) color, all 4 edges for each cell in A3:E7
$fullBlock = $ws.Range("A3:E7")
foreach ($edge in 5,6,7,8,9,10,11,12) {
  $fullBlock.Borders.Item($edge).LineStyle = 1
  $fullBlock.Borders.Item($edge).Weight = 2
  $fullBlock.Borders.Item($edge).ColorIndex = 1
}

# A3:A7 merged "Album" label re-centered (overrides block-wide left/top)
$albumBlock.HorizontalAlignment = -4108  # xlCenter
$albumBlock.VerticalAlignment = -4108    # xlCenter
$albumBlock.WrapText = $true

# B column values
$ws.Range("B3").Value = "id"
$ws.Range("B4").Value = "name"
$ws.Range("B5").Value = "album_type"
$ws.Range("B6").Value = "release_date"
$ws.Range("B7").Value = "popularity"

# B3,B5,B6 vertical-center (font will be fixed below); B4,B7 stay top
$ws.Range("B3").VerticalAlignment = -4108  # xlCenter
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("B6").VerticalAlignment = -4108

# Font for B3,B5,B6: size 11 Calibri, color FF202124
$fontCells = $ws.Range("B3,B5,B6")
foreach ($fc in $fontCells.Areas) {
  $fc.Font.Name = "Calibri"
  $fc.Font.Size = 11
  $fc.Font.Color = 2367776   # BGR encoding of RGB FF202124
}

# C/D/E column content
$ws.Range("C3").Value = "is a 36 character Universally Unique Identifier that is permanently assigned to each entity in the database, i.e. artists, release groups, releases, recordings, works, labels, areas, places and URLs."
$ws.Range("E3").Value = "6YjKAkDYmlasMqYw73iB0w"

$ws.Range("C4").Value = "album title of artist track list"
$ws.Range("E4").Value = "Bitch Please II"

$ws.Range("C5").Value = "is the total of songs in the album which will be grouped based on albums consisting of one song, and albums consisting of 8 or more songs"
$ws.Range("E5").Value = "album"

$ws.Range("C6").Value = "at which an album or single is first offered for sale. A release date is sometimes called an album launch"
$ws.Range("E6").Value = 959040000000

$ws.Range("C7").Value = "is the most heard song from the album"
$ws.Range("E7").Value = 0

# Row heights
$ws.Rows.Item(3).RowHeight = 78
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 15

Write-Output "stage1 ok"
